$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price (D) and Volume (E) columns so that
# numeric-looking strings (e.g. "0.9971") are preserved as text, matching
# the original inline-string cell type, instead of being auto-converted
# to numbers by Excel's type inference.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 49/50: Algorand and Aptos swapped positions
$ws.Range("B49").Value = 'Aptos'
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'

# Price column (D) updates
$ws.Range("D2").Value = '29.352.51'
$ws.Range("D3").Value = '1.844.93'
$ws.Range("D4").Value = '0.9971'
$ws.Range("D6").Value = '0.6271'
$ws.Range("D7").Value = '0.9993'
$ws.Range("D8").Value = '0.07484'
$ws.Range("D9").Value = '0.2901'
$ws.Range("D10").Value = '24.33'
$ws.Range("D11").Value = '0.07712'
$ws.Range("D12").Value = '1.845.22'
$ws.Range("D13").Value = '5.003'
$ws.Range("D14").Value = '0.6791'
$ws.Range("D15").Value = '0.00001027'
$ws.Range("D16").Value = '82.10'
$ws.Range("D17").Value = '2.100.69'
$ws.Range("D18").Value = '6.164'
$ws.Range("D19").Value = '29.370.60'
$ws.Range("D20").Value = '228.99'
$ws.Range("D21").Value = '12.34'
$ws.Range("D22").Value = '0.9986'
$ws.Range("D23").Value = '7.470'
$ws.Range("D24").Value = '0.9984'
$ws.Range("D25").Value = '158.95'
$ws.Range("D27").Value = '8.408'
$ws.Range("D28").Value = '17.53'
$ws.Range("D29").Value = '0.06402'
$ws.Range("D30").Value = '1.379'
$ws.Range("D31").Value = '1.470'
$ws.Range("D32").Value = '4.095'
$ws.Range("D34").Value = '1.822'
$ws.Range("D36").Value = '0.6980'
$ws.Range("D37").Value = '2.576'
$ws.Range("D38").Value = '1.260.60'
$ws.Range("D39").Value = '2.828'
$ws.Range("D41").Value = '6.589'
$ws.Range("D42").Value = '0.9112'
$ws.Range("D43").Value = '0.9986'
$ws.Range("D44").Value = '2.005.77'
$ws.Range("D45").Value = '101.45'
$ws.Range("D46").Value = '66.11'
$ws.Range("D48").Value = '1.729'
$ws.Range("D49").Value = '7.076'
$ws.Range("D50").Value = '0.1175'
$ws.Range("D51").Value = '9.007'

# Volume(1h) column (E) updates
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("E5").Value = '  -0.05%  '
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  -1.88%  '
$ws.Range("E9").Value = '  -0.19%  '
$ws.Range("E10").Value = '  -1.37%  '
$ws.Range("E11").Value = '  -0.26%  '
$ws.Range("E12").Value = '  -2.28%  '
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("E14").Value = '  +0.10%  '
$ws.Range("E15").Value = '  -3.00%  '
$ws.Range("E16").Value = '  -1.29%  '
$ws.Range("E17").Value = '  -3.94%  '
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("E19").Value = '  -0.18%  '
$ws.Range("E20").Value = '  +1.04%  '
$ws.Range("E21").Value = '  +0.13%  '
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("E24").Value = '  -0.22%  '
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("E27").Value = '  +0.17%  '
$ws.Range("E28").Value = '  -0.81%  '
$ws.Range("E29").Value = '  +14.32%  '
$ws.Range("E30").Value = '  -0.19%  '
$ws.Range("E31").Value = '  +0.43%  '
$ws.Range("E32").Value = '  -0.68%  '
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("E34").Value = '  -0.76%  '
$ws.Range("E35").Value = '  -1.88%  '
$ws.Range("E36").Value = '  +0.43%  '
$ws.Range("E37").Value = '  -0.34%  '
$ws.Range("E38").Value = '  +2.60%  '
$ws.Range("E39").Value = '  +4.00%  '
$ws.Range("E40").Value = '  +1.58%  '
$ws.Range("E41").Value = '  +3.17%  '
$ws.Range("E42").Value = '  +0.78%  '
$ws.Range("E43").Value = '  -0.19%  '
$ws.Range("E44").Value = '  -18.49%  '
$ws.Range("E45").Value = '  -0.19%  '
$ws.Range("E46").Value = '  +0.33%  '
$ws.Range("E47").Value = '  +4.55%  '
$ws.Range("E48").Value = '  +2.92%  '
$ws.Range("E49").Value = '  -1.19%  '
$ws.Range("E50").Value = '  +2.78%  '
$ws.Range("E51").Value = '  +0.28%  '

# Reset style index back to Normal (0) so only the NumberFormat-driven
# text-type coercion applied above persists, without leaving a residual
# cell style assignment on the range (matches original unstyled cells).
$ws.Range("D2:E51").Style = "Normal"

